# Scheduled-runner update: refresh market/price derived figures on the
# per-job sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR). Values below come
# from the latest pricing pull; a few rows gain/lose a HQ-profit (N) or
# NQ-profit (M) cell depending on whether that recipe has a HQ/NQ split.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 892.44446
$ws.Cells.Item(32, 9).Value = 565
$ws.Cells.Item(32, 10).Value = 986
$ws.Cells.Item(32, 11).Value = 565
$ws.Cells.Item(32, 12).Value = 986
$ws.Cells.Item(32, 13).Value = -239
$ws.Cells.Item(32, 14).Value = -1638
$ws.Cells.Item(40, 8).Value = 1629.95
$ws.Cells.Item(40, 9).Value = 1487.375
$ws.Cells.Item(40, 10).Value = 1725
$ws.Cells.Item(40, 11).Value = 1487.375
$ws.Cells.Item(40, 12).Value = 1725
$ws.Cells.Item(40, 13).Value = -1312.375
$ws.Cells.Item(40, 14).Value = -2075
$ws.Cells.Item(93, 8).Value = 35542.855
$ws.Cells.Item(93, 10).Value = 35542.855
$ws.Cells.Item(93, 12).Value = 35542.855
$ws.Cells.Item(93, 14).Value = -40534.855
$ws.Cells.Item(112, 8).Value = 2267.8262
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 10).Value = 2267.8262
$ws.Cells.Item(112, 11).Value = 0
$ws.Cells.Item(112, 12).Value = 6803.4786
$ws.Cells.Item(112, 13).Value = ""
$ws.Cells.Item(112, 14).Value = -9019.4786
$ws.Cells.Item(129, 8).Value = 922.56665
$ws.Cells.Item(129, 9).Value = 699
$ws.Cells.Item(129, 10).Value = 938.5357
$ws.Cells.Item(129, 11).Value = 2097
$ws.Cells.Item(129, 12).Value = 2815.6071
$ws.Cells.Item(129, 13).Value = 2903
$ws.Cells.Item(129, 14).Value = -12815.6071
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(23, 8).Value = 19000
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 19000
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = 19000
$ws.Cells.Item(23, 13).Value = ""
$ws.Cells.Item(23, 14).Value = -19518
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 2315006.2
$ws.Cells.Item(80, 9).Value = 18518636
$ws.Cells.Item(80, 10).Value = 201.92857
$ws.Cells.Item(80, 11).Value = 18518636
$ws.Cells.Item(80, 12).Value = 201.92857
$ws.Cells.Item(80, 13).Value = -18517638
$ws.Cells.Item(80, 14).Value = -2197.92857
$ws.Cells.Item(83, 8).Value = 2315006.2
$ws.Cells.Item(83, 9).Value = 18518636
$ws.Cells.Item(83, 10).Value = 201.92857
$ws.Cells.Item(83, 11).Value = 92593180
$ws.Cells.Item(83, 12).Value = 1009.64285
$ws.Cells.Item(83, 13).Value = -92588188
$ws.Cells.Item(83, 14).Value = -10993.64285
$ws.Cells.Item(95, 8).Value = 20780.715
$ws.Cells.Item(95, 10).Value = 20780.715
$ws.Cells.Item(95, 12).Value = 20780.715
$ws.Cells.Item(95, 14).Value = -26272.715
$ws.Cells.Item(99, 8).Value = 2165.125
$ws.Cells.Item(99, 9).Value = 1859.8
$ws.Cells.Item(99, 10).Value = 2674
$ws.Cells.Item(99, 11).Value = 1859.8
$ws.Cells.Item(99, 12).Value = 2674
$ws.Cells.Item(99, 13).Value = -361.8
$ws.Cells.Item(99, 14).Value = -5670
$ws.Cells.Item(135, 8).Value = 110040
$ws.Cells.Item(135, 10).Value = 110040
$ws.Cells.Item(135, 12).Value = 110040
$ws.Cells.Item(135, 14).Value = -120180
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(59, 8).Value = 40000
$ws.Cells.Item(59, 9).Value = 40000
$ws.Cells.Item(59, 10).Value = 40000
$ws.Cells.Item(59, 11).Value = 40000
$ws.Cells.Item(59, 12).Value = 40000
$ws.Cells.Item(59, 13).Value = -38855
$ws.Cells.Item(59, 14).Value = -42290
$ws.Cells.Item(92, 8).Value = 39565.934
$ws.Cells.Item(92, 10).Value = 39565.934
$ws.Cells.Item(92, 12).Value = 39565.934
$ws.Cells.Item(92, 14).Value = -44557.934
$ws.Cells.Item(134, 8).Value = 1647
$ws.Cells.Item(134, 9).Value = 1365.0555
$ws.Cells.Item(134, 10).Value = 2154.5
$ws.Cells.Item(134, 11).Value = 4095.1665
$ws.Cells.Item(134, 12).Value = 6463.5
$ws.Cells.Item(134, 13).Value = -1560.1665
$ws.Cells.Item(134, 14).Value = -11533.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 258
$ws.Cells.Item(17, 9).Value = 174
$ws.Cells.Item(17, 11).Value = 522
$ws.Cells.Item(17, 13).Value = -353
$ws.Cells.Item(86, 8).Value = 298
$ws.Cells.Item(86, 9).Value = 298
$ws.Cells.Item(86, 11).Value = 894
$ws.Cells.Item(86, 13).Value = 292
$ws.Cells.Item(89, 8).Value = 298
$ws.Cells.Item(89, 9).Value = 298
$ws.Cells.Item(89, 11).Value = 2682
$ws.Cells.Item(89, 13).Value = 3246
$ws.Cells.Item(113, 8).Value = 714933.3
$ws.Cells.Item(113, 10).Value = 707.8
$ws.Cells.Item(113, 12).Value = 2123.4
$ws.Cells.Item(113, 14).Value = -6463.4
$ws.Cells.Item(114, 8).Value = 1965.8334
$ws.Cells.Item(114, 9).Value = 545.6
$ws.Cells.Item(114, 10).Value = 2980.2856
$ws.Cells.Item(114, 11).Value = 1636.8
$ws.Cells.Item(114, 12).Value = 8940.856800000001
$ws.Cells.Item(114, 13).Value = 1617.2
$ws.Cells.Item(114, 14).Value = -15448.8568
$ws.Cells.Item(120, 8).Value = 6740.5625
$ws.Cells.Item(120, 9).Value = 4529.0835
$ws.Cells.Item(120, 10).Value = 13375
$ws.Cells.Item(120, 11).Value = 13587.2505
$ws.Cells.Item(120, 12).Value = 40125
$ws.Cells.Item(120, 13).Value = -8749.250499999998
$ws.Cells.Item(120, 14).Value = -49801
$ws.Cells.Item(131, 8).Value = 1549.8823
$ws.Cells.Item(131, 10).Value = 1628
$ws.Cells.Item(131, 12).Value = 4884
$ws.Cells.Item(131, 14).Value = -14964
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 853.8461
$ws.Cells.Item(22, 9).Value = 377.77777
$ws.Cells.Item(22, 11).Value = 377.77777
$ws.Cells.Item(22, 13).Value = -82.77776999999998
$ws.Cells.Item(27, 8).Value = 853.8461
$ws.Cells.Item(27, 9).Value = 377.77777
$ws.Cells.Item(27, 11).Value = 377.77777
$ws.Cells.Item(27, 13).Value = -270.77777
$ws.Cells.Item(45, 8).Value = 36523
$ws.Cells.Item(45, 10).Value = 23046
$ws.Cells.Item(45, 12).Value = 23046
$ws.Cells.Item(45, 14).Value = -23860
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(29, 8).Value = 26950
$ws.Cells.Item(29, 9).Value = 26950
$ws.Cells.Item(29, 11).Value = 26950
$ws.Cells.Item(29, 13).Value = -26660
$ws.Cells.Item(62, 8).Value = 4333.8335
$ws.Cells.Item(62, 9).Value = 4000
$ws.Cells.Item(62, 10).Value = 4667.6665
$ws.Cells.Item(62, 11).Value = 4000
$ws.Cells.Item(62, 12).Value = 4667.6665
$ws.Cells.Item(62, 13).Value = -3376
$ws.Cells.Item(62, 14).Value = -5915.6665
$ws.Cells.Item(65, 8).Value = 4333.8335
$ws.Cells.Item(65, 9).Value = 4000
$ws.Cells.Item(65, 10).Value = 4667.6665
$ws.Cells.Item(65, 11).Value = 20000
$ws.Cells.Item(65, 12).Value = 23338.3325
$ws.Cells.Item(65, 13).Value = -16880
$ws.Cells.Item(65, 14).Value = -29578.3325
$ws.Cells.Item(81, 8).Value = 100964.5
$ws.Cells.Item(81, 9).Value = 167480.17
$ws.Cells.Item(81, 10).Value = 1191
$ws.Cells.Item(81, 11).Value = 334960.34
$ws.Cells.Item(81, 12).Value = 2382
$ws.Cells.Item(81, 13).Value = -333899.34
$ws.Cells.Item(81, 14).Value = -4504
$ws.Cells.Item(84, 8).Value = 100964.5
$ws.Cells.Item(84, 9).Value = 167480.17
$ws.Cells.Item(84, 10).Value = 1191
$ws.Cells.Item(84, 11).Value = 1674801.7
$ws.Cells.Item(84, 12).Value = 11910
$ws.Cells.Item(84, 13).Value = -1669497.7
$ws.Cells.Item(84, 14).Value = -22518
$ws.Cells.Item(141, 8).Value = 57082.5
$ws.Cells.Item(141, 10).Value = 64443.332
$ws.Cells.Item(141, 12).Value = 64443.332
$ws.Cells.Item(141, 14).Value = -74803.33199999999
